$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New pre-order minifigure sets appended after the existing data (rows 86-112).
# Columns: A=Product, B=Price, C=Category, D=Image

$data = @(
    @("Naruto Set (Any 8 Minifigure)", 1600, "naruto set.jpg"),
    @("Naruto Set v2 (Any 8 Minifigure)", 1600, "naruto set 2.jpg"),
    @("Akatsuki Set (12 minifigure)", 2400, "akatsuki set.jpg"),
    @("Attack on Titan Set (8 minifigure)", 2100, "aot set.jpg"),
    @("Attack on Titan Set v2 (8 minifigure)", 2200, "aot set 2.jpg"),
    @("Bleach Set (8 minifigure)", 2150, "bleach set.jpg"),
    @("Jojo's Bizzare Adventure Set (8 minifigure)", 2000, "jojo set.jpg"),
    @("Jojo's Bizzare Adventure Set v2 (8 minifigure)", 2000, "jojo set 2.jpg"),
    @("Kaiju no.08 Set (8 minifigure)", 2050, "kaiju set.jpg"),
    @("The Boys Set (8 minifigure)", 1600, "boys set.jpg"),
    @("The Boys Set v2 (7 minifigure)", 1400, "boys set 2.jpg"),
    @("Haikyuu Set (9 minifigure)", 1600, "haikyuu set.jpg"),
    @("Loki TV Series Set (8 minifigure)", 1600, "loki set.jpg"),
    @("Loki TV Series Set v2 (8 minifigure)", 1600, "loki set 2.jpg"),
    @("Demon Slayer Set (Any 8 Minifigure)", 2000, "ds set.jpg"),
    @("Tinage Mutant Ninja Turtles Set (8 minifigure)", 1600, "tmnt set.jpg"),
    @("Football Stars Set (Any 8 Minifigure)", 1600, "football set.jpg"),
    @("Marvel's Spiderman Game Set (Any 8 Minifigure)", 1600, "spiderman set.jpg"),
    @("Spiderman-Into the Spider Verse Set (8 minifigure)", 1600, "spiderman set 2.jpg"),
    @("Deadpool & Wolverine Set (Any 8 Minifigure)", 1600, "dp wol set.jpg"),
    @("Multiverse of Madness + The Boys Set (Any 8 Minifigure)", 1600, "ran set.jpg"),
    @("X-97 Set (Any 8 Minifigure)", 1600, "xmen spiderman set.jpg"),
    @("Punisher & Daredevil Set (8 minifigure)", 1600, "punisher set.jpg"),
    @("One Piece DY Set (8 minifigure)", 2000, "one piece dy set.jpg"),
    @("One Piece TP Set (8 minifigure)", 1800, "one piece tp set.jpg"),
    @("One Piece KDL Set (Any 8 Minifigure)", 1600, "one piece kdl set.jpg"),
    @("Marvel TV series Set (Any 8 Minifigure)", 1600, "ran marvel set.jpg")
)

$startRow = 86
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $item = $data[$i]
    $ws.Cells.Item($r, 1).Value2 = $item[0]
    $ws.Cells.Item($r, 2).Value2 = $item[1]
    $ws.Cells.Item($r, 3).Value2 = "Minifigure set"
    $ws.Cells.Item($r, 4).Value2 = $item[2]
}

# Widen column A to fit the longer product names now present.
$ws.Columns.Item(1).ColumnWidth = 51.14

# Move the view / selection to reflect the newly added bottom of the list.
$ws.Range("D112").Select()
$excel.ActiveWindow.ScrollRow = 103
